$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2035928143712575
$ws.Range("C2").Value = 0.5508982035928144
$ws.Range("J2").Value = 0.01796407185628742
$ws.Range("P2").Value = 0.1467065868263473
$ws.Range("S2").Value = 0.08083832335329341
$ws.Range("B3").Value = 0.01538461538461539
$ws.Range("C3").Value = 0.02564102564102564
$ws.Range("J3").Value = 0.03076923076923077
$ws.Range("P3").Value = 0.7641025641025641
$ws.Range("S3").Value = 0.1641025641025641
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2708333333333333
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.08018867924528301
$ws.Range("F6").Value = 0.05660377358490566
$ws.Range("J6").Value = 0.2735849056603774
$ws.Range("O6").Value = 0.01415094339622642
$ws.Range("Q6").Value = 0.1650943396226415
$ws.Range("R6").Value = 0.07075471698113207
$ws.Range("S6").Value = 0.3396226415094339
$ws.Range("B7").Value = 0.1266666666666667
$ws.Range("D7").Value = 0.02666666666666667
$ws.Range("F7").Value = 0.05333333333333334
$ws.Range("J7").Value = 0.09333333333333334
$ws.Range("O7").Value = 0.02
$ws.Range("Q7").Value = 0.2533333333333334
$ws.Range("R7").Value = 0.07333333333333333
$ws.Range("S7").Value = 0.3533333333333333
$ws.Range("B8").Value = 0.1052631578947368
$ws.Range("D8").Value = 0.01894736842105263
$ws.Range("F8").Value = 0.05052631578947368
$ws.Range("J8").Value = 0.1178947368421053
$ws.Range("O8").Value = 0.01684210526315789
$ws.Range("Q8").Value = 0.208421052631579
$ws.Range("R8").Value = 0.1073684210526316
$ws.Range("S8").Value = 0.3747368421052631
$ws.Range("B9").Value = 0.07446808510638298
$ws.Range("D9").Value = 0.02127659574468085
$ws.Range("F9").Value = 0.06382978723404255
$ws.Range("J9").Value = 0.1329787234042553
$ws.Range("O9").Value = 0.01595744680851064
$ws.Range("Q9").Value = 0.175531914893617
$ws.Range("R9").Value = 0.1170212765957447
$ws.Range("S9").Value = 0.398936170212766
$ws.Range("B10").Value = 0.1163498098859316
$ws.Range("D10").Value = 0.02357414448669201
$ws.Range("E10").Value = 0.0007604562737642585
$ws.Range("F10").Value = 0.05855513307984791
$ws.Range("J10").Value = 0.1178707224334601
$ws.Range("O10").Value = 0.01749049429657795
$ws.Range("Q10").Value = 0.24106463878327
$ws.Range("R10").Value = 0.1102661596958175
$ws.Range("S10").Value = 0.3140684410646388
$ws.Range("G11").Value = 0.1353383458646616
$ws.Range("J11").Value = 0.112781954887218
$ws.Range("K11").Value = 0.1992481203007519
$ws.Range("L11").Value = 0.5375939849624061
$ws.Range("S11").Value = 0.01503759398496241
$ws.Range("G12").Value = 0.7241379310344828
$ws.Range("J12").Value = 0.1931034482758621
$ws.Range("K12").Value = 0.006896551724137931
$ws.Range("L12").Value = 0.03448275862068965
$ws.Range("S12").Value = 0.04137931034482759
$ws.Range("G13").Value = 0.4705882352941176
$ws.Range("J13").Value = 0.4117647058823529
$ws.Range("S13").Value = 0.1176470588235294
$ws.Range("F15").Value = 0.03765690376569038
$ws.Range("H15").Value = 0.1757322175732217
$ws.Range("I15").Value = 0.06694560669456066
$ws.Range("J15").Value = 0.3765690376569037
$ws.Range("K15").Value = 0.04602510460251046
$ws.Range("M15").Value = 0.008368200836820083
$ws.Range("O15").Value = 0.06276150627615062
$ws.Range("S15").Value = 0.2259414225941423
$ws.Range("F16").Value = 0.01357466063348416
$ws.Range("H16").Value = 0.2081447963800905
$ws.Range("I16").Value = 0.08144796380090498
$ws.Range("J16").Value = 0.3529411764705883
$ws.Range("K16").Value = 0.09954751131221719
$ws.Range("M16").Value = 0.009049773755656109
$ws.Range("O16").Value = 0.09502262443438914
$ws.Range("S16").Value = 0.1402714932126697
$ws.Range("F17").Value = 0.02099236641221374
$ws.Range("H17").Value = 0.200381679389313
$ws.Range("I17").Value = 0.1068702290076336
$ws.Range("J17").Value = 0.3950381679389313
$ws.Range("K17").Value = 0.08587786259541985
$ws.Range("M17").Value = 0.01717557251908397
$ws.Range("O17").Value = 0.06679389312977099
$ws.Range("S17").Value = 0.1068702290076336
$ws.Range("F18").Value = 0.008333333333333333
$ws.Range("H18").Value = 0.2416666666666667
$ws.Range("I18").Value = 0.04583333333333333
$ws.Range("J18").Value = 0.4333333333333333
$ws.Range("K18").Value = 0.06666666666666667
$ws.Range("M18").Value = 0.0125
$ws.Range("O18").Value = 0.08333333333333333
$ws.Range("S18").Value = 0.1083333333333333
$ws.Range("F19").Value = 0.02188552188552189
$ws.Range("H19").Value = 0.1902356902356902
$ws.Range("I19").Value = 0.07491582491582492
$ws.Range("J19").Value = 0.3821548821548821
$ws.Range("K19").Value = 0.09932659932659933
$ws.Range("M19").Value = 0.01683501683501683
$ws.Range("O19").Value = 0.06902356902356903
$ws.Range("S19").Value = 0.1456228956228956
